$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fertility Rate (per year) by 2010 table (rows 136:142, cols B:G) is now
# derived live from the "before 1990" table (rows 115:121) multiplied by
# 0.65 (was previously a hard-coded 0.5x snapshot).
$srcRows = @{136=115; 137=116; 138=117; 139=118; 140=119; 141=120; 142=121}
$cols = @("B","C","D","E","F","G")

foreach ($tgtRow in $srcRows.Keys) {
    $srcRow = $srcRows[$tgtRow]
    foreach ($col in $cols) {
        $ws.Range("$col$tgtRow").Formula = "=$col$srcRow*0.65"
    }
}
